# Apply the workbook edits:
#   1. Delete the "Desarquivamentos Pendentes" sheet entirely.
#   2. Rename "Paineis DARQ" -> "PAINEIS DARQ".
#   3. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO".

$wb = $excel.ActiveWorkbook

# Delete the obsolete sheet first (suppress the "delete sheet" confirmation).
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining sheets to their new (upper-case) titles.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the first sheet as the active one, as it was originally.
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
